# chore: update Sheets via scheduled runner
#
# Refreshes cached Universalis market-board pricing columns
# (currentAveragePrice / *NQ / *HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# for the leve rows whose prices moved since the last sync.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3681.3142
$ws.Range("I74").Value = 3637.6365
$ws.Range("J74").Value = 4402
$ws.Range("K74").Value = 3637.6365
$ws.Range("L74").Value = 4402
$ws.Range("M74").Value = -2701.6365
$ws.Range("N74").Value = -6274

# Row 77
$ws.Range("H77").Value = 3681.3142
$ws.Range("I77").Value = 3637.6365
$ws.Range("J77").Value = 4402
$ws.Range("K77").Value = 18188.1825
$ws.Range("L77").Value = 22010
$ws.Range("M77").Value = -13508.1825
$ws.Range("N77").Value = -31370

# Row 123
$ws.Range("H123").Value = 97016.5
$ws.Range("J123").Value = 97016.5
$ws.Range("L123").Value = 97016.5
$ws.Range("N123").Value = -106816.5

# Row 137
$ws.Range("H137").Value = 1194.8937
$ws.Range("I137").Value = 777.8889
$ws.Range("J137").Value = 1453.7241
$ws.Range("K137").Value = 2333.6667
$ws.Range("L137").Value = 4361.1723
$ws.Range("M137").Value = 216.3332999999998
$ws.Range("N137").Value = -9461.1723

# Row 138
$ws.Range("H138").Value = 1806.8586
$ws.Range("I138").Value = 597.3143
$ws.Range("J138").Value = 2468.3281
$ws.Range("K138").Value = 1791.9429
$ws.Range("L138").Value = 7404.9843
$ws.Range("M138").Value = 3348.0571
$ws.Range("N138").Value = -17684.9843

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 30355
$ws.Range("J24").Value = 30355
$ws.Range("L24").Value = 30355
$ws.Range("N24").Value = -31103

# Row 32
$ws.Range("H32").Value = 3550.4814
$ws.Range("I32").Value = 2481.5715
$ws.Range("J32").Value = 7291.6665
$ws.Range("K32").Value = 2481.5715
$ws.Range("L32").Value = 7291.6665
$ws.Range("M32").Value = -2194.5715
$ws.Range("N32").Value = -7865.6665

# Row 34
$ws.Range("H34").Value = 15622.4

# Row 100
$ws.Range("H100").Value = 30355
$ws.Range("J100").Value = 30355
$ws.Range("L100").Value = 30355
$ws.Range("N100").Value = -32519

# Row 110
$ws.Range("H110").Value = 869.1786
$ws.Range("I110").Value = 882.04346
$ws.Range("J110").Value = 810
$ws.Range("K110").Value = 882.04346
$ws.Range("L110").Value = 810
$ws.Range("M110").Value = 1162.95654
$ws.Range("N110").Value = -4900

# Row 123
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

# Row 132
$ws.Range("H132").Value = 5429.143
$ws.Range("I132").Value = 8000
$ws.Range("J132").Value = 5000.6665
$ws.Range("K132").Value = 24000
$ws.Range("L132").Value = 15001.9995
$ws.Range("M132").Value = -21470
$ws.Range("N132").Value = -20061.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1419.7333
$ws.Range("I20").Value = 1090.7222
$ws.Range("J20").Value = 1913.25
$ws.Range("K20").Value = 1090.7222
$ws.Range("L20").Value = 1913.25
$ws.Range("M20").Value = -843.7221999999999
$ws.Range("N20").Value = -2407.25

# Row 86
$ws.Range("H86").Value = 16915.285
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 16915.285
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 16915.285
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -19161.285

# Row 89
$ws.Range("H89").Value = 16915.285
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 16915.285
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 84576.425
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -95808.425

# Row 94
$ws.Range("H94").Value = 772.6786
$ws.Range("I94").Value = 782.5454999999999
$ws.Range("J94").Value = 736.5
$ws.Range("K94").Value = 782.5454999999999
$ws.Range("L94").Value = 736.5
$ws.Range("M94").Value = -331.5454999999999
$ws.Range("N94").Value = -1638.5

$ws = $wb.Worksheets.Item("CRP")
# Row 60
$ws.Range("H60").Value = 14760
$ws.Range("I60").Value = 4666.6665
$ws.Range("K60").Value = 4666.6665
$ws.Range("M60").Value = -4155.6665

# Row 105
$ws.Range("H105").Value = 388.3684
$ws.Range("I105").Value = 311.06668
$ws.Range("J105").Value = 678.25
$ws.Range("K105").Value = 311.06668
$ws.Range("L105").Value = 678.25
$ws.Range("M105").Value = 1435.93332
$ws.Range("N105").Value = -4172.25

# Row 107
$ws.Range("H107").Value = 551.7917
$ws.Range("I107").Value = 445.2
$ws.Range("K107").Value = 445.2
$ws.Range("M107").Value = 1474.8

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 55555670
$ws.Range("I2").Value = 122.76923
$ws.Range("J2").Value = 200000100
$ws.Range("K2").Value = 736.61538
$ws.Range("L2").Value = 1200000600
$ws.Range("M2").Value = -623.61538
$ws.Range("N2").Value = -1200000826

# Row 38
$ws.Range("H38").Value = 135.52
$ws.Range("I38").Value = 188.9
$ws.Range("J38").Value = 99.933334
$ws.Range("K38").Value = 566.7
$ws.Range("L38").Value = 299.800002
$ws.Range("M38").Value = -219.7
$ws.Range("N38").Value = -993.8000019999999

# Row 55
$ws.Range("H55").Value = 5950
$ws.Range("J55").Value = 5950
$ws.Range("L55").Value = 17850
$ws.Range("N55").Value = -18204

# Row 68
$ws.Range("H68").Value = 916.73
$ws.Range("I68").Value = 751.2461499999999
$ws.Range("J68").Value = 1224.0571
$ws.Range("K68").Value = 2253.73845
$ws.Range("L68").Value = 3672.1713
$ws.Range("M68").Value = -1442.73845
$ws.Range("N68").Value = -5294.1713

# Row 71
$ws.Range("H71").Value = 916.73
$ws.Range("I71").Value = 751.2461499999999
$ws.Range("J71").Value = 1224.0571
$ws.Range("K71").Value = 6761.215349999999
$ws.Range("L71").Value = 11016.5139
$ws.Range("M71").Value = -2705.215349999999
$ws.Range("N71").Value = -19128.5139

# Row 121
$ws.Range("H121").Value = 405.66666
$ws.Range("I121").Value = 163.6
$ws.Range("J121").Value = 708.25
$ws.Range("K121").Value = 490.8
$ws.Range("L121").Value = 2124.75
$ws.Range("M121").Value = 819.2
$ws.Range("N121").Value = -4744.75

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2599.889
$ws.Range("I80").Value = 2600
$ws.Range("J80").Value = 2599.8
$ws.Range("K80").Value = 2600
$ws.Range("L80").Value = 2599.8
$ws.Range("M80").Value = -1602
$ws.Range("N80").Value = -4595.8

# Row 83
$ws.Range("H83").Value = 2599.889
$ws.Range("I83").Value = 2600
$ws.Range("J83").Value = 2599.8
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 12999
$ws.Range("M83").Value = -8008
$ws.Range("N83").Value = -22983

# Row 113
$ws.Range("H113").Value = 1821.3636
$ws.Range("I113").Value = 1629.8572
$ws.Range("J113").Value = 2156.5
$ws.Range("K113").Value = 1629.8572
$ws.Range("L113").Value = 2156.5
$ws.Range("M113").Value = 540.1428000000001
$ws.Range("N113").Value = -6496.5

# Row 122
$ws.Range("H122").Value = 1235578.6
$ws.Range("I122").Value = 2222913.5
$ws.Range("K122").Value = 6668740.5
$ws.Range("M122").Value = -6666290.5

# Row 132
$ws.Range("H132").Value = 4729.143
$ws.Range("I132").Value = 4574
$ws.Range("J132").Value = 4845.5
$ws.Range("K132").Value = 13722
$ws.Range("L132").Value = 14536.5
$ws.Range("M132").Value = -11192
$ws.Range("N132").Value = -19596.5

# Row 139
$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3318.75
$ws.Range("I122").Value = 2033.3334
$ws.Range("K122").Value = 6100.0002
$ws.Range("M122").Value = -3650.0002

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 85669.664
$ws.Range("I81").Value = 144722.28
$ws.Range("J81").Value = 2996
$ws.Range("K81").Value = 289444.56
$ws.Range("L81").Value = 5992
$ws.Range("M81").Value = -288383.56
$ws.Range("N81").Value = -8114

# Row 84
$ws.Range("H84").Value = 85669.664
$ws.Range("I84").Value = 144722.28
$ws.Range("J84").Value = 2996
$ws.Range("K84").Value = 1447222.8
$ws.Range("L84").Value = 29960
$ws.Range("M84").Value = -1441918.8
$ws.Range("N84").Value = -40568

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132
$ws.Range("H132").Value = 19232884
$ws.Range("I132").Value = 23811314
$ws.Range("J132").Value = 3474.8
$ws.Range("K132").Value = 71433942
$ws.Range("L132").Value = 10424.4
$ws.Range("M132").Value = -71431412
$ws.Range("N132").Value = -15484.4
